# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# with latest values, as scraped on Thu Mar 16 23:53:12 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.016.84"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "'1.675.52"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'329.43"
$ws.Range("E5").Value = "  +7.07%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "'0.3656"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").Value = "'47.39"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.3249"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'1.151"
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("D11").Value = "'0.07295"
$ws.Range("E11").Value = "  +4.32%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'6.098"
$ws.Range("E13").Value = "  +3.60%  "
$ws.Range("D14").Value = "'19.73"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "'1.675.43"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "'6.672"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "'0.00001055"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "'0.06550"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'0.9985"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "'79.10"
$ws.Range("E20").Value = "  +3.31%  "
$ws.Range("D21").Value = "'15.88"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").Value = "'5.928"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'12.86"
$ws.Range("E23").Value = "  +3.41%  "
$ws.Range("D24").Value = "'25.027.74"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").Value = "'2.445"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "'2.402"
$ws.Range("D27").Value = "'149.17"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("D28").Value = "'18.83"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "'1.862.93"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "'126.30"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "'1.195"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("D32").Value = "'4.088"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("D33").Value = "'5.834"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("D34").Value = "'0.08480"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "'1.672"
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").Value = "'12.40"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").Value = "'5.183"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.06110"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("D39").Value = "'1.234"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("D40").Value = "'0.02242"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").Value = "'0.2090"
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("D42").Value = "'8.293"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").Value = "'0.9979"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'0.5997"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").Value = "'13.62"
$ws.Range("E45").Value = "  +7.33%  "
$ws.Range("D46").Value = "'3.836"
$ws.Range("D47").Value = "'0.5751"
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("D48").Value = "'124.52"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").Value = "'1.969"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("D50").Value = "'0.07019"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").Value = "'1.189"
$ws.Range("E51").Value = "  +3.25%  "
